$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.2633862938780678, 1.3055995132460398, -0.4116335118654511, 3.2500925253932986)"
$ws.Range("C2").Value = "NCT(3.0727906953150157, 1.9534930774721397, -0.012756880561947882, 4.1879476374518525)"
$ws.Range("D2").Value = "NCT(2.623695190221915, 1.367946381656541, -0.4800362476305731, 2.5880815160904005)"
$ws.Range("E2").Value = "NCT(3.441036129047641, 2.335691427564126, -3.1752759583429517, 4.723169862907151)"
